# Regenerate merged AHB files
# 1. Rename the "_old" / "_new" header-suffix columns to "_FV2404" / "_FV2410"
# 2. Turn the data range into an Excel Table (ListObject)
# 3. Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels -----------------------------------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}
# Column K ("diff") is unchanged.

# --- 2. Create the Excel Table over A1:U72 --------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# --- 3. Freeze panes at row 1 ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
